$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tuesday")

$ws.Cells.Item(8, 1).Value = "Trump"
$ws.Cells.Item(8, 2).Value = "Tiffany"
$ws.Cells.Item(8, 3).Value = "rnc.trumptiffany.txt"
$ws.Cells.Item(8, 4).Value = "Tuesday"
$ws.Cells.Item(8, 5).Value = "speech"

$ws.Range("A9").Select()
$ws.Activate()
